$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 73, shifting the existing rows 73-76 down to 74-77.
$ws.Rows.Item(73).Insert()

# Copy the formatting of the (now shifted) row below into the new blank row,
# so the new row matches the style of the rest of the data block (e.g. date style in column D).
$ws.Range("A74:R74").Copy()
$ws.Range("A73:R73").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row 73 with the new data record.
$ws.Cells.Item(73, 1).Value = 2
$ws.Cells.Item(73, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(73, 3).Value = "Coquimbo"
$ws.Cells.Item(73, 4).Value = 44826
$ws.Cells.Item(73, 5).Value = 4
$ws.Cells.Item(73, 6).Value = 100112022
$ws.Cells.Item(73, 7).Value = "Arveja Verde"
$ws.Cells.Item(73, 8).Value = "Perfection"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 520
$ws.Cells.Item(73, 11).Value = 28000
$ws.Cells.Item(73, 12).Value = 30000
$ws.Cells.Item(73, 13).Value = 29000
$ws.Cells.Item(73, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(73, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(73, 16).Value = 1160
$ws.Cells.Item(73, 17).Value = 25
$ws.Cells.Item(73, 18).Value = "Hortaliza"
